$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before row 284; this pushes the existing
# rows 284-330 down to 287-333 (and the sheet dimension grows to R333).
$ws.Range("A284:R286").Insert()

# Common column values shared by every data row in this block.
$mercadoId   = 9
$mercado     = "Vega Central Mapocho de Santiago"
$region      = "Metropolitana"
$codreg      = 13
$categoriaId = 100112009
$categoria   = "Acelga"
$variedad    = "Sin especificar"
$unidad      = "`$/docena de atados"
$origen      = "Región Metropolitana"
$kgUnidades  = 3
$clasif      = "Hortaliza"

# Data for the 3 newly-inserted rows (284-286).
$newRows = @(
    @{ Row=284; Fecha=44474; Calidad="Extra";   Volumen=25; PMin=12000; PMax=13000; PProm=12480; PKg=4160 },
    @{ Row=285; Fecha=44474; Calidad="Primera"; Volumen=61; PMin=10000; PMax=11000; PProm=10492; PKg=3497 },
    @{ Row=286; Fecha=44474; Calidad="Segunda"; Volumen=43; PMin=8000;  PMax=9000;  PProm=8488;  PKg=2829 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $r.Fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $categoriaId
    $ws.Cells.Item($row, 7).Value  = $categoria
    $ws.Cells.Item($row, 8).Value  = $variedad
    $ws.Cells.Item($row, 9).Value  = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Volumen
    $ws.Cells.Item($row, 11).Value = $r.PMin
    $ws.Cells.Item($row, 12).Value = $r.PMax
    $ws.Cells.Item($row, 13).Value = $r.PProm
    $ws.Cells.Item($row, 14).Value = $unidad
    $ws.Cells.Item($row, 15).Value = $origen
    $ws.Cells.Item($row, 16).Value = $r.PKg
    $ws.Cells.Item($row, 17).Value = $kgUnidades
    $ws.Cells.Item($row, 18).Value = $clasif
}
